# Insert a new data row at row 79 (pushing the existing rows 79-126 down
# to 80-127) and populate it with the new "Tercera" quality record dated
# 44873. This matches the weekly refresh of the Fruta / Hortaliza sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 79..126 down to 80..127, leaving a blank row 79 behind.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new record's data.
$ws.Range("A79").Value = 1
$ws.Range("B79").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C79").Value = "Arica y Parinacota"
$ws.Range("D79").Value = 44873
$ws.Range("E79").Value = 15
$ws.Range("F79").Value = 100112008
$ws.Range("G79").Value = "Coliflor"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Tercera"
$ws.Range("J79").Value = 1000
$ws.Range("K79").Value = 600
$ws.Range("L79").Value = 700
$ws.Range("M79").Value = 650
$ws.Range("N79").Value = "$/unidad"
$ws.Range("O79").Value = "Región de Arica y Parinacota"
$ws.Range("P79").Value = 650
$ws.Range("Q79").Value = 1
$ws.Range("R79").Value = "Hortaliza"
